# feat: add 2022-Q4 data
#
# Before: 总计 | 2021-Q4
# After:  总计 | 2022-Q4 (new data) | 2021-Q4 (old data, moved to its own new sheet)
#
# The former "2021-Q4" sheet's fund-holding data becomes the content of a brand
# new "2021-Q4" sheet, while the original sheet slot is repurposed to hold the
# new "2022-Q4" figures. The "总计" (totals) sheet gains a row for the new
# quarter and keeps the old quarter's totals in a second row.

$wb = $excel.ActiveWorkbook

$sTotal = $wb.Worksheets.Item(1)   # "总计"
$sQtr   = $wb.Worksheets.Item(2)   # currently "2021-Q4"

# ------------------------------------------------------------------
# 1. Duplicate the existing quarter sheet so its current ("2021-Q4")
#    fund-holding data is preserved on its own sheet, placed right
#    after the source sheet.
# ------------------------------------------------------------------
$sQtr.Copy($null, $sQtr)
$sOldQtr = $wb.Worksheets.Item(3)

# Free up the "2021-Q4" name on the original sheet before claiming it
# for the duplicate, then rename both sheets into their final slots.
$sQtr.Name = "2022-Q4"
$sOldQtr.Name = "2021-Q4"

# ------------------------------------------------------------------
# 2. Overwrite the (now) "2022-Q4" sheet's single data row with the
#    new quarter's fund-holding figures. Columns B/D/E/F/G hold
#    text-like values (some of which look numeric, e.g. a fund code
#    with a leading zero, or "0.21"), so force them to Text before
#    assigning to keep them as strings instead of being coerced to
#    numbers.
# ------------------------------------------------------------------
$sQtr.Range("B2").NumberFormat = "@"
$sQtr.Range("B2").Value = "001068"
$sQtr.Range("C2").Value = "国新国证新锐灵活配置混合"
$sQtr.Range("D2").NumberFormat = "@"
$sQtr.Range("D2").Value = "0.21"
$sQtr.Range("E2").NumberFormat = "@"
$sQtr.Range("E2").Value = "75.37"
$sQtr.Range("F2").NumberFormat = "@"
$sQtr.Range("F2").Value = "3.88"
$sQtr.Range("G2").NumberFormat = "@"
$sQtr.Range("G2").Value = "0.0081"
$sQtr.Range("H2").Value = 7

# ------------------------------------------------------------------
# 3. Update the "总计" sheet: the existing data row now describes the
#    new quarter, and a new row is appended describing the previous
#    quarter (using the totals that used to live in the single row).
# ------------------------------------------------------------------
$sTotal.Range("B2").Value = "2022-Q4"
$sTotal.Range("D2").Value = 0.01

# Copy A2's formatting (border/font/alignment) onto the new A3 cell
# before filling in its values, so row 3 matches row 2's look.
$sTotal.Range("A2").Copy()
$sTotal.Range("A3").PasteSpecial(-4122)  # xlPasteFormats

$sTotal.Range("A3").Value = 1
$sTotal.Range("B3").Value = "2021-Q4"
$sTotal.Range("C3").Value = 1
$sTotal.Range("D3").Value = 1.6

# ------------------------------------------------------------------
# 4. Restore "总计" as the active sheet (matches original workbook
#    view state).
# ------------------------------------------------------------------
$sTotal.Activate()
